# Update the "Funding" table on Sheet3 with refreshed grant info (h-index source data)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Select()

# Header: role column label
$ws.Range("G1").Value = "role"

# Row 2 - R21NR017312 / Tonya Palermo / sleep deficiency paper (ongoing)
$ws.Range("C2").Value = "09/21/2017-08/31/2020"
$ws.Range("D2").Value = "R21NR017312"
$ws.Range("E2").Value = "Tonya Palermo"
$ws.Range("F2").Value = "Role of sleep deficiency in self-management of pediatric chronic pain"
$ws.Range("G2").Value = "Collaborator"
$ws.Range("H2").Value = "The objective of this application is to characterize how sleep deficiency influences youths" + [char]0x2019 + " ability to engage with, implement, and benefit from a pain self-management intervention."

# Row 3 - U01CE002880 / Emily Kroshus / One Team (complete)
$ws.Range("C3").Value = "09/30/2017-09/29/2021"
$ws.Range("D3").Value = "U01CE002880"
$ws.Range("E3").Value = "Emily Kroshus"
$ws.Range("F3").Value = "One Team: Changing the culture of youth sport"
$ws.Range("G3").Value = "Collaborator"
$ws.Range("H3").Value = "Sports-related concussion in youth is increasingly being recognized as a public health concern. We propose an intervention that will utilize safety huddles to shift the culture of youth sport. At the completion of this research, we will have an intervention that will shift the culture of safety in youth sport and that can reach all youth sport stakeholders, including those in low resource communities."

# Row 4 - R21NR017312 / Tonya Palermo / sleep deficiency paper (complete)
$ws.Range("C4").Value = "09/21/2017-08/31/2020"
$ws.Range("D4").Value = "R21NR017312"
$ws.Range("E4").Value = "Tonya Palermo"
$ws.Range("F4").Value = "Role of sleep deficiency in self-management of pediatric chronic pain"
$ws.Range("G4").Value = "Collaborator"
$ws.Range("H4").Value = "The objective of this application is to characterize how sleep deficiency influences youths" + [char]0x2019 + " ability to engage with, implement, and benefit from a pain self-management intervention."

$ws.Range("A2").Select()

$wb.Save()
